# Prepend "Design: " to the answer text of six feedback rows in the table.
#
# We find each target paragraph by its exact current text (rather than doing
# a blind Find/Replace across the whole document) because a naive substring
# search for "Mostly" would also match inside the unrelated answer
# "Mostly ok." and corrupt it.
#
# Word's Paragraph.Range.Text includes trailing control characters
# (paragraph mark \r, and a cell-end mark \a for the last paragraph in a
# table cell), so those are trimmed before comparing. A couple of the
# target paragraphs contain a manual line break (<w:br/>), which Word
# renders as a literal vertical-tab character (\v) inside Range.Text
# between the two lines of the paragraph - in that case the paragraph text
# is "<target><break char><second line>", so we match either an exact
# match, or the text starting with "<target><break char>".

$d = $word.ActiveDocument

$targets = @(
    "Mostly",
    "Nothing special to mention",
    "91xx Went well, some minor budget challenges",
    "Mostly ok.",
    "Internal communication ok. External communication with suppliers mostly ok.",
    "Some things went to correct direction but regarding TK the opposite way."
)

$prefix = "Design: "
$trimChars = [char[]]@([char]13, [char]7)
$breakChar = [char]11

foreach ($wanted in $targets) {
    foreach ($p in $d.Paragraphs) {
        $text = $p.Range.Text.TrimEnd($trimChars)
        if (($text -eq $wanted) -or ($text.StartsWith($wanted + $breakChar))) {
            $p.Range.InsertBefore($prefix)
            break
        }
    }
}
